$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G7").Value = 27
$ws.Range("H7").Value = 27

$ws.Range("G6").Select()
